$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the file paths to the new OneDrive-based location
$ws.Range("B2").Value = "C:\Users\pvanausdeln\OneDrive - Blume Global\UiPath\OceanCarrierRPA\OceanCarrierRPA\APL\ContainerInformation"
$ws.Range("B3").Value = "C:\Users\pvanausdeln\OneDrive - Blume Global\UiPath\OceanCarrierRPA\OceanCarrierRPA\CMACGM\ContainerInformation"
$ws.Range("B4").Value = "C:\Users\pvanausdeln\OneDrive - Blume Global\UiPath\OceanCarrierRPA\OceanCarrierRPA\Evergreen\ContainerInformation"
$ws.Range("B5").Value = "C:\Users\pvanausdeln\OneDrive - Blume Global\UiPath\OceanCarrierRPA\OceanCarrierRPA\OOCL\ContainerInformation"

# Best-fit the column widths to the new (longer) content, matching the
# widths Excel computed for this data (A ~= 16.14 chars, B ~= 61.43 chars)
$ws.Columns.Item(1).ColumnWidth = 15.333
$ws.Columns.Item(2).ColumnWidth = 60.6665
